# "Generate Report for Handoff"
# Updates the status / handoff datetime / error-detail columns for the
# e1f29d49-765f-4abe-8a5d-7a268dab63cd file across the Overview, zh-cn and
# de-de sheets, and widens the "Error Detail" column on the locale sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78459990ababa072644d766c063513f7ef8462db/e2e/e1f29d49-765f-4abe-8a5d-7a268dab63cd.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d633c773abddb10bbb7b05fcb2da7568b9353513/e2e/e1f29d49-765f-4abe-8a5d-7a268dab63cd.md."

# ---- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 07:49:30"

# ---- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-06 07:49:19"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1").ColumnWidth = 39.17

# ---- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-06 07:49:30"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1").ColumnWidth = 39.17
